$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "All.global.sex"
$ws.Range("C1").Value = "Males.global.sex"
$ws.Range("D1").Value = "Females.global.sex"
$ws.Range("E1").Value = "Not known / missing.global.sex"
